$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + 1h volume %) per the source diff.
# A leading apostrophe forces text entry so numeric-looking strings (e.g.
# "6.70", "430.30") keep their trailing zeros instead of being coerced to numbers.

$ws.Cells.Item(2, 4).Value = "'60.866.31"  # D2
$ws.Cells.Item(2, 5).Value = "'  +0.49%  "  # E2
$ws.Cells.Item(3, 4).Value = "'2.918.18"  # D3
$ws.Cells.Item(3, 5).Value = "'  +0.77%  "  # E3
$ws.Cells.Item(4, 5).Value = "'  +0.07%  "  # E4
$ws.Cells.Item(5, 4).Value = "'593.69"  # D5
$ws.Cells.Item(5, 5).Value = "'  +1.53%  "  # E5
$ws.Cells.Item(6, 4).Value = "'145.64"  # D6
$ws.Cells.Item(6, 5).Value = "'  -0.77%  "  # E6
$ws.Cells.Item(8, 4).Value = "'0.506"  # D8
$ws.Cells.Item(8, 5).Value = "'  +0.79%  "  # E8
$ws.Cells.Item(9, 4).Value = "'6.81"  # D9
$ws.Cells.Item(9, 5).Value = "'  +2.12%  "  # E9
$ws.Cells.Item(10, 5).Value = "'  +0.01%  "  # E10
$ws.Cells.Item(11, 5).Value = "'  -2.02%  "  # E11
$ws.Cells.Item(12, 5).Value = "'  +0.30%  "  # E12
$ws.Cells.Item(13, 4).Value = "'33.62"  # D13
$ws.Cells.Item(13, 5).Value = "'  -0.83%  "  # E13
$ws.Cells.Item(14, 5).Value = "'  -0.15%  "  # E14
$ws.Cells.Item(15, 4).Value = "'3.401.70"  # D15
$ws.Cells.Item(15, 5).Value = "'  +0.77%  "  # E15
$ws.Cells.Item(16, 4).Value = "'60.891.15"  # D16
$ws.Cells.Item(16, 5).Value = "'  +0.60%  "  # E16
$ws.Cells.Item(17, 4).Value = "'6.70"  # D17
$ws.Cells.Item(17, 5).Value = "'  -1.37%  "  # E17
$ws.Cells.Item(18, 4).Value = "'2.921.50"  # D18
$ws.Cells.Item(19, 4).Value = "'430.30"  # D19
$ws.Cells.Item(19, 5).Value = "'  +1.46%  "  # E19
$ws.Cells.Item(20, 4).Value = "'13.36"  # D20
$ws.Cells.Item(20, 5).Value = "'  -1.93%  "  # E20
$ws.Cells.Item(21, 4).Value = "'0.681"  # D21
$ws.Cells.Item(21, 5).Value = "'  +1.62%  "  # E21
$ws.Cells.Item(22, 4).Value = "'7.06"  # D22
$ws.Cells.Item(22, 5).Value = "'  -0.15%  "  # E22
$ws.Cells.Item(23, 4).Value = "'81.55"  # D23
$ws.Cells.Item(23, 5).Value = "'  +1.64%  "  # E23
$ws.Cells.Item(24, 4).Value = "'10.97"  # D24
$ws.Cells.Item(24, 5).Value = "'  -0.68%  "  # E24
$ws.Cells.Item(25, 5).Value = "'  -0.68%  "  # E25
$ws.Cells.Item(26, 4).Value = "'11.92"  # D26
$ws.Cells.Item(26, 5).Value = "'  +0.64%  "  # E26
$ws.Cells.Item(27, 5).Value = "'  -0.03%  "  # E27
$ws.Cells.Item(28, 4).Value = "'2.30"  # D28
$ws.Cells.Item(28, 5).Value = "'  +4.33%  "  # E28
$ws.Cells.Item(29, 5).Value = "'  +0.12%  "  # E29
$ws.Cells.Item(30, 5).Value = "'  +0.03%  "  # E30
$ws.Cells.Item(31, 4).Value = "'7.04"  # D31
$ws.Cells.Item(31, 5).Value = "'  -3.50%  "  # E31
$ws.Cells.Item(32, 4).Value = "'26.38"  # D32
$ws.Cells.Item(32, 5).Value = "'  -0.01%  "  # E32
$ws.Cells.Item(33, 5).Value = "'  +0.30%  "  # E33
$ws.Cells.Item(34, 4).Value = "'0.0₃0849"  # D34
$ws.Cells.Item(34, 5).Value = "'  +2.14%  "  # E34
$ws.Cells.Item(35, 5).Value = "'  +0.61%  "  # E35
$ws.Cells.Item(36, 4).Value = "'5.61"  # D36
$ws.Cells.Item(36, 5).Value = "'  -1.06%  "  # E36
$ws.Cells.Item(37, 4).Value = "'3.02"  # D37
$ws.Cells.Item(37, 5).Value = "'  +2.81%  "  # E37
$ws.Cells.Item(38, 5).Value = "'  -0.27%  "  # E38
$ws.Cells.Item(39, 5).Value = "'  -1.80%  "  # E39
$ws.Cells.Item(40, 4).Value = "'8.56"  # D40
$ws.Cells.Item(40, 5).Value = "'  -1.89%  "  # E40
$ws.Cells.Item(41, 5).Value = "'  -2.17%  "  # E41
$ws.Cells.Item(42, 4).Value = "'40.31"  # D42
$ws.Cells.Item(42, 5).Value = "'  -3.12%  "  # E42
$ws.Cells.Item(43, 4).Value = "'373.66"  # D43
$ws.Cells.Item(43, 5).Value = "'  +0.44%  "  # E43
$ws.Cells.Item(44, 5).Value = "'  -0.17%  "  # E44
$ws.Cells.Item(45, 4).Value = "'2.697.29"  # D45
$ws.Cells.Item(45, 5).Value = "'  +1.79%  "  # E45
$ws.Cells.Item(46, 4).Value = "'130.90"  # D46
$ws.Cells.Item(46, 5).Value = "'  -1.66%  "  # E46
$ws.Cells.Item(47, 5).Value = "'  -0.04%  "  # E47
$ws.Cells.Item(48, 4).Value = "'23.93"  # D48
$ws.Cells.Item(48, 5).Value = "'  -5.59%  "  # E48
$ws.Cells.Item(49, 5).Value = "'  -0.07%  "  # E49
$ws.Cells.Item(50, 5).Value = "'  -3.63%  "  # E50
$ws.Cells.Item(51, 5).Value = "'  +2.25%  "  # E51
